$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"2"
$ws.Range("F2").Value = [double]"0.6666666666666666"
$ws.Range("G2").Value = [double]"0.09286699999999999"
$ws.Range("H2").Value = [double]"0.278601"
$ws.Range("I2").Value = [double]"0.003009076821730935"
$ws.Range("J2").Value = [double]"0.003071957783644885"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"0.506715"
$ws.Range("N2").Value = [double]"1.520145"
$ws.Range("O2").Value = [double]"0.003122343715987576"
$ws.Range("P2").Value = [double]"0.003132472094339857"
$ws.Range("Q2").Value = [double]"0.047057101905"
$ws.Range("R2").Value = [double]"0.423513917145"
$ws.Range("S2").Value = [double]"9.395372105255453E-06"
$ws.Range("T2").Value = [double]"9.622822032257719E-06"
$ws.Range("E3").Value = [double]"2"
$ws.Range("F3").Value = [double]"0.6666666666666666"
$ws.Range("G3").Value = [double]"0.09286699999999999"
$ws.Range("H3").Value = [double]"0.278601"
$ws.Range("I3").Value = [double]"0.003009076821730935"
$ws.Range("J3").Value = [double]"0.003071957783644885"
$ws.Range("M3").Value = [double]"88.13219433333332"
$ws.Range("N3").Value = [double]"264.396583"
$ws.Range("O3").Value = [double]"0.5430646480820168"
$ws.Range("P3").Value = [double]"0.5448262620252092"
$ws.Range("Q3").Value = [double]"8.184572491153665"
$ws.Range("R3").Value = [double]"73.66115242038299"
$ws.Range("S3").Value = [double]"0.001634123245245064"
$ws.Range("T3").Value = [double]"0.001673683276362489"
$ws.Range("E4").Value = [double]"2"
$ws.Range("F4").Value = [double]"0.6666666666666666"
$ws.Range("G4").Value = [double]"0.09286699999999999"
$ws.Range("H4").Value = [double]"0.278601"
$ws.Range("I4").Value = [double]"0.003009076821730935"
$ws.Range("J4").Value = [double]"0.003071957783644885"
$ws.Range("M4").Value = [double]"1.5741895"
$ws.Range("N4").Value = [double]"3.148379"
$ws.Range("O4").Value = [double]"0.009700049718478087"
$ws.Range("P4").Value = [double]"0.006487676741301404"
$ws.Range("Q4").Value = [double]"0.1461902562965"
$ws.Range("R4").Value = [double]"0.8771415377789999"
$ws.Range("S4").Value = [double]"2.918819477751009E-05"
$ws.Range("T4").Value = [double]"1.992986906321273E-05"
$ws.Range("E5").Value = [double]"2"
$ws.Range("F5").Value = [double]"0.6666666666666666"
$ws.Range("G5").Value = [double]"0.09286699999999999"
$ws.Range("H5").Value = [double]"0.278601"
$ws.Range("I5").Value = [double]"0.003009076821730935"
$ws.Range("J5").Value = [double]"0.003071957783644885"
$ws.Range("M5").Value = [double]"72.07364666666666"
$ws.Range("N5").Value = [double]"216.22094"
$ws.Range("O5").Value = [double]"0.4441129584835175"
$ws.Range("P5").Value = [double]"0.4455535891391496"
$ws.Range("Q5").Value = [double]"6.693263344993333"
$ws.Range("R5").Value = [double]"60.23937010493999"
$ws.Range("S5").Value = [double]"0.001336370009603105"
$ws.Range("T5").Value = [double]"0.001368721816186926"
$ws.Range("I6").Value = [double]"0.9272539658256183"
$ws.Range("J6").Value = [double]"0.9466308793322996"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"0.506715"
$ws.Range("N6").Value = [double]"1.520145"
$ws.Range("O6").Value = [double]"0.003122343715987576"
$ws.Range("P6").Value = [double]"0.003132472094339857"
$ws.Range("Q6").Value = [double]"14.500754532605"
$ws.Range("R6").Value = [double]"130.506790793445"
$ws.Range("S6").Value = [double]"0.002895205593320178"
$ws.Range("T6").Value = [double]"0.002965294813148829"
$ws.Range("I7").Value = [double]"0.9272539658256183"
$ws.Range("J7").Value = [double]"0.9466308793322996"
$ws.Range("M7").Value = [double]"88.13219433333332"
$ws.Range("N7").Value = [double]"264.396583"
$ws.Range("O7").Value = [double]"0.5430646480820168"
$ws.Range("P7").Value = [double]"0.5448262620252092"
$ws.Range("Q7").Value = [double]"2522.094898409377"
$ws.Range("R7").Value = [double]"22698.8540856844"
$ws.Range("S7").Value = [double]"0.5035588486337438"
$ws.Range("T7").Value = [double]"0.5157493635042537"
$ws.Range("I8").Value = [double]"0.9272539658256183"
$ws.Range("J8").Value = [double]"0.9466308793322996"
$ws.Range("M8").Value = [double]"1.5741895"
$ws.Range("N8").Value = [double]"3.148379"
$ws.Range("O8").Value = [double]"0.009700049718478087"
$ws.Range("P8").Value = [double]"0.006487676741301404"
$ws.Range("Q8").Value = [double]"45.04886480033983"
$ws.Range("R8").Value = [double]"270.293188802039"
$ws.Range("S8").Value = [double]"0.008994409570164478"
$ws.Range("T8").Value = [double]"0.006141435138441857"
$ws.Range("I9").Value = [double]"0.9272539658256183"
$ws.Range("J9").Value = [double]"0.9466308793322996"
$ws.Range("M9").Value = [double]"72.07364666666666"
$ws.Range("N9").Value = [double]"216.22094"
$ws.Range("O9").Value = [double]"0.4441129584835175"
$ws.Range("P9").Value = [double]"0.4455535891391496"
$ws.Range("Q9").Value = [double]"2062.544543940949"
$ws.Range("R9").Value = [double]"18562.90089546854"
$ws.Range("S9").Value = [double]"0.4118055020283897"
$ws.Range("T9").Value = [double]"0.4217747858764553"
$ws.Range("E10").Value = [double]"3"
$ws.Range("F10").Value = [double]"1"
$ws.Range("G10").Value = [double]"0.243062"
$ws.Range("H10").Value = [double]"0.729186"
$ws.Range("I10").Value = [double]"0.007875695677081898"
$ws.Range("J10").Value = [double]"0.00804027483183793"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"0.506715"
$ws.Range("N10").Value = [double]"1.520145"
$ws.Range("O10").Value = [double]"0.003122343715987576"
$ws.Range("P10").Value = [double]"0.003132472094339857"
$ws.Range("Q10").Value = [double]"0.12316316133"
$ws.Range("R10").Value = [double]"1.10846845197"
$ws.Range("S10").Value = [double]"2.459062890636718E-05"
$ws.Range("T10").Value = [double]"2.51859365415554E-05"
$ws.Range("E11").Value = [double]"3"
$ws.Range("F11").Value = [double]"1"
$ws.Range("G11").Value = [double]"0.243062"
$ws.Range("H11").Value = [double]"0.729186"
$ws.Range("I11").Value = [double]"0.007875695677081898"
$ws.Range("J11").Value = [double]"0.00804027483183793"
$ws.Range("M11").Value = [double]"88.13219433333332"
$ws.Range("N11").Value = [double]"264.396583"
$ws.Range("O11").Value = [double]"0.5430646480820168"
$ws.Range("P11").Value = [double]"0.5448262620252092"
$ws.Range("Q11").Value = [double]"21.42158741904866"
$ws.Range("R11").Value = [double]"192.794286771438"
$ws.Range("S11").Value = [double]"0.004277011901275542"
$ws.Range("T11").Value = [double]"0.004380552882285628"
$ws.Range("E12").Value = [double]"3"
$ws.Range("F12").Value = [double]"1"
$ws.Range("G12").Value = [double]"0.243062"
$ws.Range("H12").Value = [double]"0.729186"
$ws.Range("I12").Value = [double]"0.007875695677081898"
$ws.Range("J12").Value = [double]"0.00804027483183793"
$ws.Range("M12").Value = [double]"1.5741895"
$ws.Range("N12").Value = [double]"3.148379"
$ws.Range("O12").Value = [double]"0.009700049718478087"
$ws.Range("P12").Value = [double]"0.006487676741301404"
$ws.Range("Q12").Value = [double]"0.382625648249"
$ws.Range("R12").Value = [double]"2.295753889494"
$ws.Range("S12").Value = [double]"7.639463963529735E-05"
$ws.Range("T12").Value = [double]"5.2162704020186E-05"
$ws.Range("E13").Value = [double]"3"
$ws.Range("F13").Value = [double]"1"
$ws.Range("G13").Value = [double]"0.243062"
$ws.Range("H13").Value = [double]"0.729186"
$ws.Range("I13").Value = [double]"0.007875695677081898"
$ws.Range("J13").Value = [double]"0.00804027483183793"
$ws.Range("M13").Value = [double]"72.07364666666666"
$ws.Range("N13").Value = [double]"216.22094"
$ws.Range("O13").Value = [double]"0.4441129584835175"
$ws.Range("P13").Value = [double]"0.4455535891391496"
$ws.Range("Q13").Value = [double]"17.51836470609333"
$ws.Range("R13").Value = [double]"157.66528235484"
$ws.Range("S13").Value = [double]"0.003497698507264691"
$ws.Range("T13").Value = [double]"0.003582373308990562"
$ws.Range("G14").Value = [double]"1.8951925"
$ws.Range("H14").Value = [double]"3.790385"
$ws.Range("I14").Value = [double]"0.06140803325689756"
$ws.Range("J14").Value = [double]"0.04179418847656979"
$ws.Range("K14").Value = [double]"3"
$ws.Range("L14").Value = [double]"1"
$ws.Range("M14").Value = [double]"0.506715"
$ws.Range("N14").Value = [double]"1.520145"
$ws.Range("O14").Value = [double]"0.003122343715987576"
$ws.Range("P14").Value = [double]"0.003132472094339857"
$ws.Range("Q14").Value = [double]"0.9603224676375001"
$ws.Range("R14").Value = [double]"5.761934805825001"
$ws.Range("S14").Value = [double]"0.0001917369867508302"
$ws.Range("T14").Value = [double]"0.0001309191291084353"
$ws.Range("G15").Value = [double]"1.8951925"
$ws.Range("H15").Value = [double]"3.790385"
$ws.Range("I15").Value = [double]"0.06140803325689756"
$ws.Range("J15").Value = [double]"0.04179418847656979"
$ws.Range("M15").Value = [double]"88.13219433333332"
$ws.Range("N15").Value = [double]"264.396583"
$ws.Range("O15").Value = [double]"0.5430646480820168"
$ws.Range("P15").Value = [double]"0.5448262620252092"
$ws.Range("Q15").Value = [double]"167.0274737090758"
$ws.Range("R15").Value = [double]"1002.164842254455"
$ws.Range("S15").Value = [double]"0.03334853197006586"
$ws.Range("T15").Value = [double]"0.02277057148206659"
$ws.Range("G16").Value = [double]"1.8951925"
$ws.Range("H16").Value = [double]"3.790385"
$ws.Range("I16").Value = [double]"0.06140803325689756"
$ws.Range("J16").Value = [double]"0.04179418847656979"
$ws.Range("M16").Value = [double]"1.5741895"
$ws.Range("N16").Value = [double]"3.148379"
$ws.Range("O16").Value = [double]"0.009700049718478087"
$ws.Range("P16").Value = [double]"0.006487676741301404"
$ws.Range("Q16").Value = [double]"2.98339213397875"
$ws.Range("R16").Value = [double]"11.933568535915"
$ws.Range("S16").Value = [double]"0.0005956609757058622"
$ws.Range("T16").Value = [double]"0.000271147184501009"
$ws.Range("G17").Value = [double]"1.8951925"
$ws.Range("H17").Value = [double]"3.790385"
$ws.Range("I17").Value = [double]"0.06140803325689756"
$ws.Range("J17").Value = [double]"0.04179418847656979"
$ws.Range("M17").Value = [double]"72.07364666666666"
$ws.Range("N17").Value = [double]"216.22094"
$ws.Range("O17").Value = [double]"0.4441129584835175"
$ws.Range("P17").Value = [double]"0.4455535891391496"
$ws.Range("Q17").Value = [double]"136.5934346103166"
$ws.Range("R17").Value = [double]"819.5606076619"
$ws.Range("S17").Value = [double]"0.027272103324375"
$ws.Range("T17").Value = [double]"0.01862155068089376"
$ws.Range("E18").Value = [double]"1"
$ws.Range("F18").Value = [double]"0.3333333333333333"
$ws.Range("G18").Value = [double]"0.01398766666666667"
$ws.Range("H18").Value = [double]"0.041963"
$ws.Range("I18").Value = [double]"0.000453228418671488"
$ws.Range("J18").Value = [double]"0.0004626995756479349"
$ws.Range("K18").Value = [double]"3"
$ws.Range("L18").Value = [double]"1"
$ws.Range("M18").Value = [double]"0.506715"
$ws.Range("N18").Value = [double]"1.520145"
$ws.Range("O18").Value = [double]"0.003122343715987576"
$ws.Range("P18").Value = [double]"0.003132472094339857"
$ws.Range("Q18").Value = [double]"0.007087760515"
$ws.Range("R18").Value = [double]"0.06378984463500001"
$ws.Range("S18").Value = [double]"1.415134904945907E-06"
$ws.Range("T18").Value = [double]"1.44939350878005E-06"
$ws.Range("E19").Value = [double]"1"
$ws.Range("F19").Value = [double]"0.3333333333333333"
$ws.Range("G19").Value = [double]"0.01398766666666667"
$ws.Range("H19").Value = [double]"0.041963"
$ws.Range("I19").Value = [double]"0.000453228418671488"
$ws.Range("J19").Value = [double]"0.0004626995756479349"
$ws.Range("M19").Value = [double]"88.13219433333332"
$ws.Range("N19").Value = [double]"264.396583"
$ws.Range("O19").Value = [double]"0.5430646480820168"
$ws.Range("P19").Value = [double]"0.5448262620252092"
$ws.Range("Q19").Value = [double]"1.232763756936555"
$ws.Range("R19").Value = [double]"11.094873812429"
$ws.Range("S19").Value = [double]"0.0002461323316866006"
$ws.Range("T19").Value = [double]"0.0002520908802409149"
$ws.Range("E20").Value = [double]"1"
$ws.Range("F20").Value = [double]"0.3333333333333333"
$ws.Range("G20").Value = [double]"0.01398766666666667"
$ws.Range("H20").Value = [double]"0.041963"
$ws.Range("I20").Value = [double]"0.000453228418671488"
$ws.Range("J20").Value = [double]"0.0004626995756479349"
$ws.Range("M20").Value = [double]"1.5741895"
$ws.Range("N20").Value = [double]"3.148379"
$ws.Range("O20").Value = [double]"0.009700049718478087"
$ws.Range("P20").Value = [double]"0.006487676741301404"
$ws.Range("Q20").Value = [double]"0.02201923799616667"
$ws.Range("R20").Value = [double]"0.132115427977"
$ws.Range("S20").Value = [double]"4.396338194940636E-06"
$ws.Range("T20").Value = [double]"3.001845275141137E-06"
$ws.Range("E21").Value = [double]"1"
$ws.Range("F21").Value = [double]"0.3333333333333333"
$ws.Range("G21").Value = [double]"0.01398766666666667"
$ws.Range("H21").Value = [double]"0.041963"
$ws.Range("I21").Value = [double]"0.000453228418671488"
$ws.Range("J21").Value = [double]"0.0004626995756479349"
$ws.Range("M21").Value = [double]"72.07364666666666"
$ws.Range("N21").Value = [double]"216.22094"
$ws.Range("O21").Value = [double]"0.4441129584835175"
$ws.Range("P21").Value = [double]"0.4455535891391496"
$ws.Range("Q21").Value = [double]"1.008142145024444"
$ws.Range("R21").Value = [double]"9.07327930522"
$ws.Range("S21").Value = [double]"0.0002012846138850008"
$ws.Range("T21").Value = [double]"0.0002061574566230988"

Write-Output "applied updates"
